$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated crime statistics data (rows 2-14, columns B:G)

# Row 2
$ws.Range("B2").Value = 38
$ws.Range("C2").Value = 228
$ws.Range("D2").Value = 85
$ws.Range("E2").Value = 133

# Row 3
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 146
$ws.Range("D3").Value = 55
$ws.Range("E3").Value = 68

# Row 4
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 170
$ws.Range("D4").Value = 76
$ws.Range("G4").Value = 3

# Row 5
$ws.Range("B5").Value = 82
$ws.Range("C5").Value = 298
$ws.Range("D5").Value = 133
$ws.Range("E5").Value = 175
$ws.Range("G5").Value = 1

# Row 6
$ws.Range("B6").Value = 50
$ws.Range("C6").Value = 154
$ws.Range("D6").Value = 74
$ws.Range("E6").Value = 81
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 6

# Row 7
$ws.Range("B7").Value = 25
$ws.Range("C7").Value = 119
$ws.Range("D7").Value = 32
$ws.Range("E7").Value = 68
$ws.Range("F7").Value = 7

# Row 8
$ws.Range("B8").Value = 25
$ws.Range("C8").Value = 140
$ws.Range("D8").Value = 45
$ws.Range("E8").Value = 73
$ws.Range("F8").Value = 6

# Row 9
$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 102
$ws.Range("D9").Value = 49
$ws.Range("E9").Value = 61
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 3

# Row 10
$ws.Range("B10").Value = 38
$ws.Range("C10").Value = 203
$ws.Range("D10").Value = 64
$ws.Range("E10").Value = 91
$ws.Range("G10").Value = 7

# Row 11
$ws.Range("B11").Value = 29
$ws.Range("C11").Value = 91

# Row 12
$ws.Range("B12").Value = 51
$ws.Range("C12").Value = 248
$ws.Range("E12").Value = 141

# Row 13
$ws.Range("C13").Value = 377
$ws.Range("D13").Value = 139
$ws.Range("E13").Value = 177
$ws.Range("F13").Value = 18

# Row 14
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = 60
$ws.Range("D14").Value = 17

# Update the selection to match the final cursor position used by the author
$ws.Range("I19").Select()

$wb.Save()
